$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-TextFromBase64($b64) {
    $bytes = [System.Convert]::FromBase64String($b64)
    return [System.Text.Encoding]::UTF8.GetString($bytes)
}

# C3: fix one au_orcid entry (NA -> ORCID) for "search_author" bug fix
$c3b64 = "bGlzdChhdV9pZCA9IGMoImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDExODkwNzU2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDMyNDY2MzUzIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDg1OTE3ODU5IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDQ1OTAxNjUyIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDQ1MDg4Njc4IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDgzMjg1NTUwIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDExOTEzNjc0IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDA2OTE2MjE1IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDA1MzQ3MjQ3IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDEzNzEzNDEzIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDQ5OTM0MzM4IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDI1ODMwMjU2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDYxMjE2MTI5IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDgyODQ3MjM2IiwgCiJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA1MDU4MDY0NCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAyODMzMTM0MyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA2NjU2NDA5MCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAyNDczMzc2NiIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA3MjExMTcyMSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAzMDUzMzgyMiIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAzOTk1ODIyNyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA2NzYyODA4OSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAwNzA5MDc1OSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA4MzI4NTU1MCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAwNDk1NjMzNiIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA3NzI4ODAwNyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA2Njk2NTcwNSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA1NjE1OTAxNiIsIAoiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwMjk3NTkwNjUiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwNTk0MDc0NjEiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwMTE4OTA3NTYiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwNjMwMzkyNTUiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwMTcxNjg5NTciLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwNzg4OTM0OTMiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwNTY4NDE2MTYiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwNTcyMzM4MzAiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwMDA5NDY2MjkiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwMzA1MzM4MjIiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwMjU4MzAyNTYiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwODYwNzIxOTYiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwMTE5MTM2NzQiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvQTUwODgxNzc1MDQiLCAKImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDI0ODI2MTEyIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDQzNTcyNDA5IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDkxNTA3NTM5IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDU4NzQyMDQ1IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDI2NDIxNDkxIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDM5NDA5NTE1IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDYwOTEwMjUwIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDIxNTY0NDUyIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDg0NjkzNzA3IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDYxMTkwNzYzIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDc4MTk0OTE3IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDYxNjE1MDMzIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDkxNTAxNjg1IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0E1MDc1NDU3NDU0IiwgCiJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAyMjkxMDM2MyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA3MTY2Mzk0MiIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAyMjE4MzY5MyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA2MjEyODE5MCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA4ODQ2NjI3NCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA0NTEwMTQ0NyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA3MjUxNTk4NSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA2NjA3MjA5OCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA4MjI1NzYzNSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTAwOTg4NTM2OSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA3NjM2MDg3MCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9BNTA3MzMxMTk3NSIpLCBhdV9kaXNwbGF5X25hbWUgPSBjKCJDaHJpc3RpbmUgTS4gRHVyYW5kIiwgIk5pbmEgTWFydGluZXoiLCAiS2FybCBOZXVtYW5uIiwgCiJSb2JlcnQgQmVuZWRpY3QiLCAiQXJ0aHVyIFcuIEJha2VyIiwgIkNhbWVyb24gUi4gV29sZmUiLCAiVmFsZW50aW5hIFN0b3NvciIsICJBbmVlc2hhIFNoZXR0eSIsICJaYWNoYXJ5IEMuIERpZXRjaCIsICJMZWFoIEdvdWR5IiwgIk1pY2hlbGxlIEEuIENhbGxlZ2FyaSIsICJBbGxhbiBCLiBNYXNzaWUiLCAiRGlhbmUgQnJvd24iLCAiV2lsbGEgQ29jaHJhbiIsICJBYmltZXJla2kgRC4gTXV6YWFsZSIsICJEZXJlayBNLiBGaW5lIiwgIkFhcm9uIEEuUi4gVG9iaWFuIiwgIkNoZXJ5bCBBLiBXaW5rbGVyIiwgIkZhd2F6IEFsIEFtbWFyeSIsICJEb3JyeSBMLiBTZWdldiIsICJOZWhhIEFncmF3YWwiLCAiTWFyY3VzIFIuIFBlcmVpcmEiLCAiS2FydGhpayBSYW5nYW5uYSIsICJDYW1lcm9uIFIuIFdvbGZlIiwgIlJhY2hlbCBGcmllZG1hbuKAkE1vcmFjbyIsICJXaWxsaWFtIEguIEtpdGNoZW5zIiwgIk9sdXdhZmlzYXlvIEFkZWJpeWkiLCAiQ2hhbmRyYXNoZWtoYXIgQS4gS3ViYWwiLCAKIkFuZHJldyBNLiBDYW1lcm9uIiwgIk5pcmFqIE0uIERlc2FpIiwgIkNocmlzdGluZSBNLiBEdXJhbmQiLCAiU2hhbmUgT3R0bWFubiIsICJOYWhlbCBFbGlhcyIsICJBbGV4YW5kZXIgR2lsYmVydCIsICJDb2xlbWFuIEkuIFNtaXRoIiwgIkpvc2UgQS4gQ2FzdGlsbG8tTHVnbyIsICJTYW5kZXIgRmxvcm1hbiIsICJEb3JyeSBMLiBTZWdldiIsICJBbGxhbiBCLiBNYXNzaWUiLCAiU2FwbmEgQS4gTWVodGEiLCAiVmFsZW50aW5hIFN0b3NvciIsICJKb25hdGhhbiBIYW5kIiwgIkVtaWx5IEEuIEJsdW1iZXJnIiwgIkNhcmxvcyBTYW50b3MiLCAiUm9jaGVsbGUgR29sZGJlcmciLCAiU2hpa2hhIE1laHRhIiwgIlJvYmVydCBNLiBDYW5ub24iLCAiRW1tYW5vdWlsIEdpb3JnYWtpcyIsICJKb2FubmEgU2NoYWVubWFuIiwgIlNhaW1hIEFzbGFtIiwgIlBldGVyIEcuIFN0b2NrIiwgIkplbm5pZmVyIFByaWNlIiwgIlNlbnUgQXBld29raW4iLCAiRXN0aGVyIEJlbmFtdSIsICJNYXJpbyBTcGFnZ2lhcmkiLCAKIkpvaG4gVy4gQmFkZGxleSIsICJNaWNoZWxlIE1vcnJpcyIsICJKYWNxdWVzIFNpbWtpbnMiLCAiVGltb3RoeSBMLiBQcnVldHQiLCAiR2hhZHkgSGFpZGFyIiwgIkRhdmlkIFdvamNpZWNob3dza2kiLCAiQXZpbmFzaCBLdW1hciDDgWdhcndhbCIsICJWYXNhbnRoaSBCYWxhcmFtYW4iLCAiR2F1cmF2IEd1cHRhIiwgIldpbGwgQ2hhcG1hbiIsICJUaGFuZ2FtYW5pIE11dGh1a3VtYXIiLCAiQ2F0aGVyaW5lIEIuIFNtYWxsIiwgIk1hcmljYXIgTWFsaW5pcyIpLCBhdV9vcmNpZCA9IGMoImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMy0yNjA1LTkyNTciLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAzLTQ2ODYtMDY3NCIsIE5BLCBOQSwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi0wOTE0LTAyOTEiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTUzNjUtNTAzMCIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMDM4OC01NDYwIiwgTkEsIE5BLCBOQSwgCk5BLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTUyODgtNTEyNSIsIE5BLCBOQSwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi0zMjg3LTYwNjEiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAzLTI4MzgtMjA1OCIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMDUxNy0zNzY2IiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMS01NTUyLTA5MTciLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAzLTA5MjgtMjg2MCIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMTkyNC00ODAxIiwgTkEsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItOTIxNC05MTIyIiwgTkEsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItNTM2NS01MDMwIiwgTkEsIE5BLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTMyMTItNTgzMSIsIE5BLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTU4MTAtMjM5OCIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMDQ1Ny0zNTA2IiwgCiJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDMtMjYwNS05MjU3IiwgTkEsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDEtNjQ2Ni03MzQ3IiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMS01MDY5LTE4ODAiLCBOQSwgTkEsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMTYzNS05MTM2IiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi0xOTI0LTQ4MDEiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTUyODgtNTEyNSIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItNTU4OC05MDVYIiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi0wMzg4LTU0NjAiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTU3NTItOTU3NiIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItNTE5My02MTcwIiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi02ODc0LTY3MzYiLCBOQSwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMS02NjE2LTQ4NjYiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTAzNDUtMDY0MyIsIAoiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTUwMTktNTQ5NyIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDMtMTE3NC0zOTYxIiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi00MDUxLTU2MjkiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTU4MDYtMDE2NyIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMzgyNS01Njk2IiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi0zOTE1LTU4MTQiLCBOQSwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMy0zMDAwLTYwOTQiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAxLTkxMTEtNjI1WCIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMzI1NS01NzI3IiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMS05NjI2LTA3NjAiLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTA3MTUtODUzNSIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDMtMDYzNC04MjExIiwgTkEsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDEtNTI4My05NjMxIiwgCk5BLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAzLTE5MTktMTk3MCIsIE5BLCAiaHR0cHM6Ly9vcmNpZC5vcmcvMDAwMC0wMDAyLTAyOTQtNjIwOSIsICJodHRwczovL29yY2lkLm9yZy8wMDAwLTAwMDItMDYwMS02NjE1IiwgImh0dHBzOi8vb3JjaWQub3JnLzAwMDAtMDAwMi01NzIwLTk5OTQiKSwgYXV0aG9yX3Bvc2l0aW9uID0gYygiZmlyc3QiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAKIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgIm1pZGRsZSIsICJtaWRkbGUiLCAibWlkZGxlIiwgImxhc3QiKSwgYXVfYWZmaWxpYXRpb25fcmF3ID0gYygiRGVwYXJ0bWVudCBvZiBNZWRpY2luZSwgSm9obnMgSG9wa2lucyBVbml2ZXJzaXR5IFNjaG9vbCBvZiBNZWRpY2luZSwgQmFsdGltb3JlLCBNRCwgVVNBIiwgIkRvbm9yIDEsIEdBLCBVU0EuIE9SQ2lEIElEOiAwMDAwLTAwMDMtNDY4Ni0wNjc0IiwgCiJEb25vciAyLCBWQSwgVVNBIiwgIkRvbm9yIDMsIElMLCBVU0EiLCAiRGl2aXNpb24gb2YgSW5mZWN0aW91cyBEaXNlYXNlcywgRGVwYXJ0bWVudCBvZiBNZWRpY2luZSwgRHVrZSBVbml2ZXJzaXR5IFNjaG9vbCBvZiBNZWRpY2luZSwgRHVyaGFtLCBOQywgVVNBIiwgIkRpdmlzaW9uIG9mIEluZmVjdGlvdXMgRGlzZWFzZXMsIERlcGFydG1lbnQgb2YgTWVkaWNpbmUsIER1a2UgVW5pdmVyc2l0eSBTY2hvb2wgb2YgTWVkaWNpbmUsIER1cmhhbSwgTkMsIFVTQSIsICJEaXZpc2lvbnMgb2YgSW5mZWN0aW91cyBEaXNlYXNlcyBhbmQgT3JnYW4gVHJhbnNwbGFudGF0aW9uLCBGZWluYmVyZyBTY2hvb2wgb2YgTWVkaWNpbmUsIE5vcnRod2VzdGVybiBVbml2ZXJzaXR5LCBDaGljYWdvLCBJTCwgVVNBIiwgIkRpdmlzaW9ucyBvZiBJbmZlY3Rpb3VzIERpc2Vhc2VzIGFuZCBPcmdhbiBUcmFuc3BsYW50YXRpb24sIEZlaW5iZXJnIFNjaG9vbCBvZiBNZWRpY2luZSwgTm9ydGh3ZXN0ZXJuIFVuaXZlcnNpdHksIENoaWNhZ28sIElMLCBVU0EiLCAKIkRpdmlzaW9ucyBvZiBJbmZlY3Rpb3VzIERpc2Vhc2VzIGFuZCBPcmdhbiBUcmFuc3BsYW50YXRpb24sIEZlaW5iZXJnIFNjaG9vbCBvZiBNZWRpY2luZSwgTm9ydGh3ZXN0ZXJuIFVuaXZlcnNpdHksIENoaWNhZ28sIElMLCBVU0EiLCAiRGl2aXNpb25zIG9mIEluZmVjdGlvdXMgRGlzZWFzZXMgYW5kIE9yZ2FuIFRyYW5zcGxhbnRhdGlvbiwgRmVpbmJlcmcgU2Nob29sIG9mIE1lZGljaW5lLCBOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSwgQ2hpY2FnbywgSUwsIFVTQSIsICJEaXZpc2lvbnMgb2YgSW5mZWN0aW91cyBEaXNlYXNlcyBhbmQgT3JnYW4gVHJhbnNwbGFudGF0aW9uLCBGZWluYmVyZyBTY2hvb2wgb2YgTWVkaWNpbmUsIE5vcnRod2VzdGVybiBVbml2ZXJzaXR5LCBDaGljYWdvLCBJTCwgVVNBIiwgIkRlcGFydG1lbnQgb2YgU3VyZ2VyeSwgTllVIEdyb3NzbWFuIFNjaG9vbCBvZiBNZWRpY2luZSwgTmV3IFlvcmssIE5ZLCBVU0EiLCAiRGVwYXJ0bWVudCBvZiBNZWRpY2luZSwgSm9obnMgSG9wa2lucyBVbml2ZXJzaXR5IFNjaG9vbCBvZiBNZWRpY2luZSwgQmFsdGltb3JlLCBNRCwgVVNBIiwgCiJEZXBhcnRtZW50IG9mIE1lZGljaW5lLCBKb2hucyBIb3BraW5zIFVuaXZlcnNpdHkgU2Nob29sIG9mIE1lZGljaW5lLCBCYWx0aW1vcmUsIE1ELCBVU0EiLCAiRGVwYXJ0bWVudCBvZiBTdXJnZXJ5LCBKb2hucyBIb3BraW5zIFVuaXZlcnNpdHkgU2Nob29sIG9mIE1lZGljaW5lLCBCYWx0aW1vcmUsIE1ELCBVU0EiLCAiRGVwYXJ0bWVudCBvZiBNZWRpY2luZSwgSm9obnMgSG9wa2lucyBVbml2ZXJzaXR5IFNjaG9vbCBvZiBNZWRpY2luZSwgQmFsdGltb3JlLCBNRCwgVVNBIiwgIkRlcGFydG1lbnQgb2YgUGF0aG9sb2d5LCBKb2hucyBIb3BraW5zIFVuaXZlcnNpdHkgU2Nob29sIG9mIE1lZGljaW5lLCBCYWx0aW1vcmUsIE1ELCBVU0EiLCAiQ2FuY2VyIElubm92YXRpb24gTGFib3JhdG9yeSwgQ2VudGVyIGZvciBDYW5jZXIgUmVzZWFyY2gsIE5DSSBhbmQgQmFzaWMgUmVzZWFyY2ggUHJvZ3JhbSwgRnJlZGVyaWNrIE5hdGlvbmFsIExhYm9yYXRvcmllcyBmb3IgQ2FuY2VyIFJlc2VhcmNoLCBGcmVkZXJpY2ssIE1ELCBVU0EiLCAKIkRlcGFydG1lbnQgb2YgTWVkaWNpbmUsIEpvaG5zIEhvcGtpbnMgVW5pdmVyc2l0eSBTY2hvb2wgb2YgTWVkaWNpbmUsIEJhbHRpbW9yZSwgTUQsIFVTQSIsICJEZXBhcnRtZW50IG9mIFN1cmdlcnksIE5ZVSBHcm9zc21hbiBTY2hvb2wgb2YgTWVkaWNpbmUsIE5ldyBZb3JrLCBOWSwgVVNBIiwgIiIsICIiLCAiIiwgIkRpdmlzaW9uIG9mIEluZmVjdGlvdXMgRGlzZWFzZXMsIERlcGFydG1lbnQgb2YgTWVkaWNpbmUsIER1a2UgVW5pdmVyc2l0eSBTY2hvb2wgb2YgTWVkaWNpbmUsIER1cmhhbSwgTkMsIFVTQSIsICIiLCAiIiwgIiIsICIiLCAiRGl2aXNpb24gb2YgSW5mZWN0aW91cyBEaXNlYXNlcywgRGVwYXJ0bWVudCBvZiBNZWRpY2luZSwgRHVrZSBVbml2ZXJzaXR5IFNjaG9vbCBvZiBNZWRpY2luZSwgRHVyaGFtLCBOQywgVVNBIiwgIiIsICJEZXBhcnRtZW50IG9mIE1lZGljaW5lLCBKb2hucyBIb3BraW5zIFVuaXZlcnNpdHkgU2Nob29sIG9mIE1lZGljaW5lLCBCYWx0aW1vcmUsIE1ELCBVU0EiLCAKIiIsICIiLCAiIiwgIiIsICIiLCAiIiwgIkRlcGFydG1lbnQgb2YgU3VyZ2VyeSwgTllVIEdyb3NzbWFuIFNjaG9vbCBvZiBNZWRpY2luZSwgTmV3IFlvcmssIE5ZLCBVU0EiLCAiRGVwYXJ0bWVudCBvZiBTdXJnZXJ5LCBOWVUgR3Jvc3NtYW4gU2Nob29sIG9mIE1lZGljaW5lLCBOZXcgWW9yaywgTlksIFVTQSIsICIiLCAiRGl2aXNpb25zIG9mIEluZmVjdGlvdXMgRGlzZWFzZXMgYW5kIE9yZ2FuIFRyYW5zcGxhbnRhdGlvbiwgRmVpbmJlcmcgU2Nob29sIG9mIE1lZGljaW5lLCBOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSwgQ2hpY2FnbywgSUwsIFVTQSIsICIiLCAiIiwgIiIsICIiLCAiIiwgIiIsICIiLCAiIiwgIiIsICIiLCAiIiwgIiIsICIiLCAiIiwgIiIsICIiLCAiIiwgIiIsICIiLCAiIiwgIiIsICIiLCAiIiwgIkRlcGFydG1lbnQgb2YgTWVkaWNpbmUsIEpvaG5zIEhvcGtpbnMgVW5pdmVyc2l0eSBTY2hvb2wgb2YgTWVkaWNpbmUsIEJhbHRpbW9yZSwgTUQsIFVTQSIsICIiLCAKIiIsICIiKSwgaW5zdGl0dXRpb25faWQgPSBjKCJodHRwczovL29wZW5hbGV4Lm9yZy9JMjc5OTg1MzQzNiIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JNDIxMDExMTUwNCIsIE5BLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTQyMTAxMTE1MDQiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTE3MDg5NzMxNyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMTcwODk3MzE3IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kxMTE5Nzk5MjEiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTExMTk3OTkyMSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMTExOTc5OTIxIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kxMTE5Nzk5MjEiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTExMTk3OTkyMSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JNTcyMDY5NzQiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTI3OTk4NTM0MzYiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTI3OTk4NTM0MzYiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTI3OTk4NTM0MzYiLCAKImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0k0MjEwMTMwNjQ5IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0k1NzIwNjk3NCIsIE5BLCBOQSwgTkEsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMTcwODk3MzE3IiwgTkEsIE5BLCBOQSwgTkEsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMTcwODk3MzE3IiwgTkEsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMjc5OTg1MzQzNiIsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsICJodHRwczovL29wZW5hbGV4Lm9yZy9JNTcyMDY5NzQiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTU3MjA2OTc0IiwgTkEsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMTExOTc5OTIxIiwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCAKTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgTkEsIE5BLCBOQSksIGluc3RpdHV0aW9uX2Rpc3BsYXlfbmFtZSA9IGMoIkpvaG5zIEhvcGtpbnMgTWVkaWNpbmUiLCAiRG9ub3IgTmV0d29yayBXZXN0IiwgTkEsICJEb25vciBOZXR3b3JrIFdlc3QiLCAiRHVrZSBVbml2ZXJzaXR5IiwgIkR1a2UgVW5pdmVyc2l0eSIsICJOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSIsICJOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSIsICJOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSIsICJOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSIsICJOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSIsICJOZXcgWW9yayBVbml2ZXJzaXR5IiwgIkpvaG5zIEhvcGtpbnMgTWVkaWNpbmUiLCAiSm9obnMgSG9wa2lucyBNZWRpY2luZSIsICJKb2hucyBIb3BraW5zIE1lZGljaW5lIiwgIkpvaG5zIEhvcGtpbnMgTWVkaWNpbmUiLCAiSm9obnMgSG9wa2lucyBNZWRpY2luZSIsIAoiRnJlZGVyaWNrIE5hdGlvbmFsIExhYm9yYXRvcnkgZm9yIENhbmNlciBSZXNlYXJjaCIsICJKb2hucyBIb3BraW5zIE1lZGljaW5lIiwgIk5ldyBZb3JrIFVuaXZlcnNpdHkiLCBOQSwgTkEsIE5BLCAiRHVrZSBVbml2ZXJzaXR5IiwgTkEsIE5BLCBOQSwgTkEsICJEdWtlIFVuaXZlcnNpdHkiLCBOQSwgIkpvaG5zIEhvcGtpbnMgTWVkaWNpbmUiLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCAiTmV3IFlvcmsgVW5pdmVyc2l0eSIsICJOZXcgWW9yayBVbml2ZXJzaXR5IiwgTkEsICJOb3J0aHdlc3Rlcm4gVW5pdmVyc2l0eSIsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgIkpvaG5zIEhvcGtpbnMgTWVkaWNpbmUiLCBOQSwgTkEsIE5BKSwgaW5zdGl0dXRpb25fcm9yID0gYygiaHR0cHM6Ly9yb3Iub3JnLzAzN3pnbjM1NCIsICJodHRwczovL3Jvci5vcmcvMDFzcnZ0ejk4IiwgCk5BLCAiaHR0cHM6Ly9yb3Iub3JnLzAxc3J2dHo5OCIsICJodHRwczovL3Jvci5vcmcvMDBweTgxNDE1IiwgImh0dHBzOi8vcm9yLm9yZy8wMHB5ODE0MTUiLCAiaHR0cHM6Ly9yb3Iub3JnLzAwMGUwYmU0NyIsICJodHRwczovL3Jvci5vcmcvMDAwZTBiZTQ3IiwgImh0dHBzOi8vcm9yLm9yZy8wMDBlMGJlNDciLCAiaHR0cHM6Ly9yb3Iub3JnLzAwMGUwYmU0NyIsICJodHRwczovL3Jvci5vcmcvMDAwZTBiZTQ3IiwgImh0dHBzOi8vcm9yLm9yZy8wMTkwYWs1NzIiLCAiaHR0cHM6Ly9yb3Iub3JnLzAzN3pnbjM1NCIsICJodHRwczovL3Jvci5vcmcvMDM3emduMzU0IiwgImh0dHBzOi8vcm9yLm9yZy8wMzd6Z24zNTQiLCAiaHR0cHM6Ly9yb3Iub3JnLzAzN3pnbjM1NCIsICJodHRwczovL3Jvci5vcmcvMDM3emduMzU0IiwgImh0dHBzOi8vcm9yLm9yZy8wM3Y2bTMyMDkiLCAiaHR0cHM6Ly9yb3Iub3JnLzAzN3pnbjM1NCIsICJodHRwczovL3Jvci5vcmcvMDE5MGFrNTcyIiwgTkEsIApOQSwgTkEsICJodHRwczovL3Jvci5vcmcvMDBweTgxNDE1IiwgTkEsIE5BLCBOQSwgTkEsICJodHRwczovL3Jvci5vcmcvMDBweTgxNDE1IiwgTkEsICJodHRwczovL3Jvci5vcmcvMDM3emduMzU0IiwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgImh0dHBzOi8vcm9yLm9yZy8wMTkwYWs1NzIiLCAiaHR0cHM6Ly9yb3Iub3JnLzAxOTBhazU3MiIsIE5BLCAiaHR0cHM6Ly9yb3Iub3JnLzAwMGUwYmU0NyIsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgImh0dHBzOi8vcm9yLm9yZy8wMzd6Z24zNTQiLCBOQSwgTkEsIE5BKSwgaW5zdGl0dXRpb25fY291bnRyeV9jb2RlID0gYygiVVMiLCAiVVMiLCBOQSwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgIlVTIiwgCiJVUyIsIE5BLCBOQSwgTkEsICJVUyIsIE5BLCBOQSwgTkEsIE5BLCAiVVMiLCBOQSwgIlVTIiwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgIlVTIiwgIlVTIiwgTkEsICJVUyIsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgIlVTIiwgTkEsIE5BLCBOQSksIGluc3RpdHV0aW9uX3R5cGUgPSBjKCJoZWFsdGhjYXJlIiwgIm5vbnByb2ZpdCIsIE5BLCAibm9ucHJvZml0IiwgImVkdWNhdGlvbiIsICJlZHVjYXRpb24iLCAiZWR1Y2F0aW9uIiwgImVkdWNhdGlvbiIsICJlZHVjYXRpb24iLCAiZWR1Y2F0aW9uIiwgImVkdWNhdGlvbiIsICJlZHVjYXRpb24iLCAiaGVhbHRoY2FyZSIsICJoZWFsdGhjYXJlIiwgImhlYWx0aGNhcmUiLCAiaGVhbHRoY2FyZSIsICJoZWFsdGhjYXJlIiwgImZhY2lsaXR5IiwgImhlYWx0aGNhcmUiLCAiZWR1Y2F0aW9uIiwgTkEsIE5BLCBOQSwgCiJlZHVjYXRpb24iLCBOQSwgTkEsIE5BLCBOQSwgImVkdWNhdGlvbiIsIE5BLCAiaGVhbHRoY2FyZSIsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsICJlZHVjYXRpb24iLCAiZWR1Y2F0aW9uIiwgTkEsICJlZHVjYXRpb24iLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsICJoZWFsdGhjYXJlIiwgTkEsIE5BLCBOQSksIGluc3RpdHV0aW9uX2xpbmVhZ2UgPSBjKCJodHRwczovL29wZW5hbGV4Lm9yZy9JMjc5OTg1MzQzNiIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JNDIxMDExMTUwNCIsIE5BLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTQyMTAxMTE1MDQiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTE3MDg5NzMxNyIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMTcwODk3MzE3IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kxMTE5Nzk5MjEiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTExMTk3OTkyMSIsIAoiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTExMTk3OTkyMSIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMTExOTc5OTIxIiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kxMTE5Nzk5MjEiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTU3MjA2OTc0IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0k0MjEwMTMwNjQ5LCBodHRwczovL29wZW5hbGV4Lm9yZy9JNDIxMDE0MDg4NCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMjc5OTg1MzQzNiIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JNTcyMDY5NzQiLCBOQSwgTkEsIE5BLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvSTE3MDg5NzMxNyIsIApOQSwgTkEsIE5BLCBOQSwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kxNzA4OTczMTciLCBOQSwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kyNzk5ODUzNDM2IiwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0k1NzIwNjk3NCIsICJodHRwczovL29wZW5hbGV4Lm9yZy9JNTcyMDY5NzQiLCBOQSwgImh0dHBzOi8vb3BlbmFsZXgub3JnL0kxMTE5Nzk5MjEiLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsIE5BLCBOQSwgTkEsICJodHRwczovL29wZW5hbGV4Lm9yZy9JMjc5OTg1MzQzNiIsIE5BLCBOQSwgTkEpKQ=="
$ws.Range("C3").Value = Get-TextFromBase64 $c3b64

# D4: updated/cleaned abstract text (exclude more columns / cleanup)
$d4b64 = "V2hpbGUgbGl2aW5nIGRvbm9yIChMRCkga2lkbmV5IHRyYW5zcGxhbnRhdGlvbiBpcyB0aGUgb3B0aW1hbCB0cmVhdG1lbnQgZm9yIHBhdGllbnRzIHdpdGgga2lkbmV5IGZhaWx1cmUsIExEcyBhc3N1bWUgYSBoaWdoZXIgcmlzayBvZiBmdXR1cmUga2lkbmV5IGZhaWx1cmUgdGhlbXNlbHZlcy4gTERzIG9mIEFmcmljYW4gYW5jZXN0cnkgaGF2ZSBhbiBldmVuIGdyZWF0ZXIgcmlzayBvZiBraWRuZXkgZmFpbHVyZSBwb3N0LWRvbmF0aW9uIHRoYW4gV2hpdGUgTERzLiBCZWNhdXNlIGV2aWRlbmNlIHN1Z2dlc3RzIHRoYXQgQXBvbGlwb3Byb3RlaW4gTDEgKEFQT0wxKSByaXNrIHZhcmlhbnRzIGNvbnRyaWJ1dGUgdG8gdGhpcyBncmVhdGVyIHJpc2ssIHRyYW5zcGxhbnQgbmVwaHJvbG9naXN0cyBhcmUgaW5jcmVhc2luZ2x5IHVzaW5nIEFQT0wxIGdlbmV0aWMgdGVzdGluZyB0byBldmFsdWF0ZSBMRCBjYW5kaWRhdGVzIG9mIEFmcmljYW4gYW5jZXN0cnkuIEhvd2V2ZXIsIG5lcGhyb2xvZ2lzdHMgZG8gbm90IGNvbnNpc3RlbnRseSBwZXJmb3JtIGdlbmV0aWMgY291bnNlbGxpbmcgd2l0aCBMRCBjYW5kaWRhdGVzIGFib3V0IEFQT0wxIGR1ZSB0byBhIGxhY2sgb2Yga25vd2xlZGdlIGFuZCBza2lsbCBpbiBjb3Vuc2VsbGluZy4gV2l0aG91dCBwcm9wZXIgY291bnNlbGxpbmcsIEFQT0wxIHRlc3Rpbmcgd2lsbCBtYWduaWZ5IExEIGNhbmRpZGF0ZXMnIGRlY2lzaW9uYWwgY29uZmxpY3QgYWJvdXQgZG9uYXRpbmcsIGplb3BhcmRpc2luZyB0aGVpciBpbmZvcm1lZCBjb25zZW50LiBHaXZlbiBjdWx0dXJhbCBjb25jZXJucyBhYm91dCBnZW5ldGljIHRlc3RpbmcgYW1vbmcgcGVvcGxlIG9mIEFmcmljYW4gYW5jZXN0cnksIHByb3RlY3RpbmcgTEQgY2FuZGlkYXRlcycgc2FmZXR5IGlzIGVzc2VudGlhbCB0byBpbXByb3ZlIGluZm9ybWVkIGRlY2lzaW9ucyBhYm91dCBkb25hdGluZy4gQ2xpbmljYWwgJ2NoYXRib3RzJywgbW9iaWxlIGFwcHMgdGhhdCBwcm92aWRlIGdlbmV0aWMgaW5mb3JtYXRpb24gdG8gcGF0aWVudHMsIGNhbiBpbXByb3ZlIGluZm9ybWVkIHRyZWF0bWVudCBkZWNpc2lvbnMuIE5vIGNoYXRib3Qgb24gQVBPTDEgaXMgYXZhaWxhYmxlIGFuZCBubyBuZXBocm9sb2dpc3QgdHJhaW5pbmcgcHJvZ3JhbW1lcyBhcmUgYXZhaWxhYmxlIHRvIHByb3ZpZGUgY3VsdHVyYWxseSBjb21wZXRlbnQgY291bnNlbGxpbmcgdG8gTERzIGFib3V0IEFQT0wxLiBHaXZlbiB0aGUgc2hvcnRhZ2Ugb2YgZ2VuZXRpYyBjb3Vuc2VsbG9ycywgaW5jcmVhc2luZyBuZXBocm9sb2dpc3RzJyBnZW5ldGljIGxpdGVyYWN5IGlzIGNyaXRpY2FsIHRvIGludGVncmF0aW5nIGdlbmV0aWMgdGVzdGluZyBpbnRvIHByYWN0aWNlLlVzaW5nIGEgbm9uLXJhbmRvbWlzZWQsIHByZS1wb3N0IHRyaWFsIGRlc2lnbiBpbiB0d28gdHJhbnNwbGFudCBjZW50cmVzIChDaGljYWdvLCBJTCwgYW5kIFdhc2hpbmd0b24sIERDKSwgd2Ugd2lsbCBldmFsdWF0ZSB0aGUgZWZmZWN0aXZlbmVzcyBvZiBjdWx0dXJhbGx5IGNvbXBldGVudCBBUE9MMSB0ZXN0aW5nLCBjaGF0Ym90IGFuZCBjb3Vuc2VsbGluZyBvbiBMRCBjYW5kaWRhdGVzJyBkZWNpc2lvbmFsIGNvbmZsaWN0IGFib3V0IGRvbmF0aW5nLCBwcmVwYXJlZG5lc3MgZm9yIGRlY2lzaW9uLW1ha2luZywgd2lsbGluZ25lc3MgdG8gZG9uYXRlIGFuZCBzYXRpc2ZhY3Rpb24gd2l0aCBpbmZvcm1lZCBjb25zZW50IGFuZCBsb25naXR1ZGluYWxseSBldmFsdWF0ZSB0aGUgaW1wbGVtZW50YXRpb24gb2YgdGhpcyBpbnRlcnZlbnRpb24gaW50byBjbGluaWNhbCBwcmFjdGljZSB1c2luZyB0aGUgUmVhY2gsIEVmZmVjdGl2ZW5lc3MsIEFkb3B0aW9uLCBJbXBsZW1lbnRhdGlvbiBhbmQgTWFpbnRlbmFuY2UgZnJhbWV3b3JrLlRoaXMgc3R1ZHkgd2lsbCBjcmVhdGUgYSBtb2RlbCBmb3IgQVBPTDEgdGVzdGluZyBvZiBMRHMgb2YgQWZyaWNhbiBhbmNlc3RyeSwgd2hpY2ggY2FuIGJlIGltcGxlbWVudGVkIG5hdGlvbmFsbHkgdmlhIGltcGxlbWVudGF0aW9uIHNjaWVuY2UgYXBwcm9hY2hlcy4gQVBPTDEgd2lsbCBzZXJ2ZSBhcyBhIG1vZGVsIGZvciBpbnRlZ3JhdGluZyBjdWx0dXJhbGx5IGNvbXBldGVudCBnZW5ldGljIHRlc3RpbmcgaW50byB0cmFuc3BsYW50IGFuZCBvdGhlciBwcmFjdGljZXMgdG8gaW1wcm92ZSBpbmZvcm1lZCBjb25zZW50LiBUaGlzIHN0dWR5IGludm9sdmVzIGh1bWFuIHBhcnRpY2lwYW50cyBhbmQgd2FzIGFwcHJvdmVkIGJ5IE5vcnRod2VzdGVybiBVbml2ZXJzaXR5IElSQiAoU1RVMDAyMTQwMzgpLiBQYXJ0aWNpcGFudHMgZ2F2ZSBpbmZvcm1lZCBjb25zZW50IHRvIHBhcnRpY2lwYXRlIGluIHRoZSBzdHVkeSBiZWZvcmUgdGFraW5nIHBhcnQuQ2xpbmljYWxUcmlhbHMuZ292IElkZW50aWZpZXI6IE5DVDA0OTEwODY3LiBSZWdpc3RlcmVkIDggTWF5IDIwMjEsIGh0dHBzOi8vcmVnaXN0ZXIuZ292L3Bycy9hcHAvYWN0aW9uL1NlbGVjdFByb3RvY29sP3NpZD1TMDAwQVdaNiZzZWxlY3RhY3Rpb249RWRpdCZ1aWQ9VTAwMDFQUEYmdHM9NyZjeD0tOGp2N20yIENsaW5pY2FsVHJpYWxzLmdvdiBJZGVudGlmaWVyOiBOQ1QwNDk5OTQzNi4gUmVnaXN0ZXJlZCA1IE5vdmVtYmVyIDIwMjEsIGh0dHBzOi8vcmVnaXN0ZXIuZ292L3Bycy9hcHAvYWN0aW9uL1NlbGVjdFByb3RvY29sP3NpZD1TMDAwQVlXVyZzZWxlY3RhY3Rpb249RWRpdCZ1aWQ9VTAwMDFQUEYmdHM9MTEmY3g9OXRueTd2Lg=="
$ws.Range("D4").Value = Get-TextFromBase64 $d4b64

# AF4: updated related_works ids
$af4b64 = "YygiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzIwNTU5MTU1ODEiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzIzNTU3Njk1MzgiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzIwODE4MTYyNTIiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzIzNzM4ODUxNjgiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzIxMTI0MzYzMDgiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzIzNjIzOTEyOTQiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzMxMzk5MTQ0OTQiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzE1NTY4MTk5MjYiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzMxNzk5NjUyNzMiLCAiaHR0cHM6Ly9vcGVuYWxleC5vcmcvVzIxNTEzMzM1NzEiKQ=="
$ws.Range("AF4").Value = Get-TextFromBase64 $af4b64
